$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells remain text (they store numeric-looking strings as text)
$ws.Range("D2:D11").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Move Robot21 to location (6, 12) and remove the toolkit."
$ws.Range("B2").Value = 201.94161
$ws.Range("C2").Value = 69560
$ws.Range("D2").Value = "0.06801"
$ws.Range("E2").Value = "9de2b1c8-62bb-4ec9-8933-49f8cf82bf59"

# Row 3
$ws.Range("A3").Value = "Move Robot41 to location (1, 11) and remove the liquid spill."
$ws.Range("B3").Value = 49.854476
$ws.Range("C3").Value = 6151
$ws.Range("D3").Value = "0.01233"
$ws.Range("E3").Value = "bf0c468c-37e6-424a-8599-b9ca644eb9fe"

# Row 4
$ws.Range("A4").Value = "Move Robot9 to location (1, 5) and remove the large debris."
$ws.Range("B4").Value = 23.678993
$ws.Range("C4").Value = 3814
$ws.Range("D4").Value = "0.0075"
$ws.Range("E4").Value = "4897631a-4f02-4640-af14-a35e73be6209"

# Row 5
$ws.Range("A5").Value = "Move Robot42 to location (1, 11) and remove the dust."
$ws.Range("B5").Value = 31.254502
$ws.Range("C5").Value = 4532
$ws.Range("D5").Value = "0.00933"
$ws.Range("E5").Value = "524813e7-5ac8-42e0-a338-7ff843e65e9b"

# Row 6
$ws.Range("A6").Value = "Move Robot32 to location (3, 9) and remove the grass."
$ws.Range("B6").Value = 115.464501
$ws.Range("C6").Value = 29142
$ws.Range("D6").Value = "0.03906"
$ws.Range("E6").Value = "362fc2fc-917b-4eba-ae03-47f4d9b5492f"

# Row 7
$ws.Range("A7").Value = "Move Robot14 to location (11, 12) and remove the small debris."
$ws.Range("B7").Value = 49.087498
$ws.Range("C7").Value = 7910
$ws.Range("D7").Value = "0.01452"
$ws.Range("E7").Value = "7c59c425-3e7a-4e80-a86a-827e062ff418"

# Row 8
$ws.Range("A8").Value = "Move Robot39 to location (6, 4) and remove the vehicle."
$ws.Range("B8").Value = 185.681
$ws.Range("C8").Value = 62879
$ws.Range("D8").Value = "0.05877"
$ws.Range("E8").Value = "b4ebd4d7-6e9c-46aa-94e6-3df823eea15a"

# Row 9
$ws.Range("A9").Value = "Move Robot15 to location (11, 2) and remove the construction materials."
$ws.Range("B9").Value = 34.808001
$ws.Range("C9").Value = 4609
$ws.Range("D9").Value = "0.0099"
$ws.Range("E9").Value = "c475182f-d963-43b4-9177-2b593667df6b"

# Row 10
$ws.Range("A10").Value = "Move Robot2 to location (1, 10) and remove the tree branches."
$ws.Range("B10").Value = 54.454471
$ws.Range("C10").Value = 6312
$ws.Range("D10").Value = "0.01305"
$ws.Range("E10").Value = "8d836cf1-428b-4fae-bdcd-3043a4118542"

# Row 11
$ws.Range("A11").Value = "Move Robot26 to location (1, 3) and remove the screws."
$ws.Range("B11").Value = 156.042001
$ws.Range("C11").Value = 44203
$ws.Range("D11").Value = "0.04989"
$ws.Range("E11").Value = "041e1fea-45e2-48f8-a406-67aa87dbfa3d"
